# Update Delivery_results sheet with new simulation results
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Delivery_results")

$ws1.Range("D2").Value = 62
$ws1.Range("E2").Value = 100
$ws1.Range("F2").Value = 1.26
$ws1.Range("G2").Value = 6.42

$ws1.Range("D3").Value = 46
$ws1.Range("E3").Value = 76.67
$ws1.Range("F3").Value = 0.38
$ws1.Range("G3").Value = 1.96

$ws1.Range("D4").Value = 41
$ws1.Range("E4").Value = 95.34999999999999
$ws1.Range("F4").Value = 1.31
$ws1.Range("G4").Value = 6.39

$ws1.Range("D5").Value = 18
$ws1.Range("E5").Value = 100
$ws1.Range("F5").Value = 1.4
$ws1.Range("G5").Value = 6.61

$ws1.Range("D6").Value = 29
$ws1.Range("E6").Value = 78.38
$ws1.Range("F6").Value = 0.27
$ws1.Range("G6").Value = 1.69

# Add a new "Total_distance" sheet right after "Total_emissions"
$ws2 = $wb.Worksheets.Item("Total_emissions")
$wsNew = $wb.Worksheets.Add($null, $ws2)
$wsNew.Name = "Total_distance"

# Match the bold/centered/bordered header style used on the other sheets
# by copying the format from an existing header cell (e.g. Delivery_results!A1)
$ws1.Range("A1").Copy()
$wsNew.Range("A1").PasteSpecial(-4122)

$wsNew.Range("A1").Value = "Total distance (km)"
$wsNew.Range("A2").Value = 39.76204
